$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boswell_2012")

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 14).Value2
}
